$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting for columns B:E of data rows so numeric-looking
# strings (prices) are not auto-converted to numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '75.853.28'
$ws.Cells.Item(2, 5).Value = '  +2.14%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.835.40'
$ws.Cells.Item(3, 5).Value = '  +7.45%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 5).Value = '  +0.05%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '191.59'
$ws.Cells.Item(5, 5).Value = '  +3.46%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '596.90'
$ws.Cells.Item(6, 5).Value = '  +2.39%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.550'
$ws.Cells.Item(8, 5).Value = '  +3.62%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.194'
$ws.Cells.Item(9, 5).Value = '  +1.72%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '2.835.41'
$ws.Cells.Item(10, 5).Value = '  +7.49%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.389'
$ws.Cells.Item(11, 5).Value = '  +10.08%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -1.98%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '4.89'
$ws.Cells.Item(13, 5).Value = '  +4.59%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '3.366.44'
$ws.Cells.Item(14, 5).Value = '  +6.81%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '75.574.05'
$ws.Cells.Item(15, 5).Value = '  +1.82%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0000189'
$ws.Cells.Item(16, 5).Value = '  +2.97%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '27.23'
$ws.Cells.Item(17, 5).Value = '  +4.02%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.843.52'
$ws.Cells.Item(18, 5).Value = '  +7.04%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -2.26%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '12.42'
$ws.Cells.Item(20, 5).Value = '  +4.87%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '380.66'
$ws.Cells.Item(21, 5).Value = '  +3.29%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +3.60%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '4.15'
$ws.Cells.Item(23, 5).Value = '  +2.52%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '71.40'
$ws.Cells.Item(24, 5).Value = '  +2.32%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.00%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '4.22'
$ws.Cells.Item(26, 5).Value = '  +3.68%  '

# Row 27
$ws.Cells.Item(27, 2).Value = 'WrappedeETH'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(27, 4).Value = '2.961.72'
$ws.Cells.Item(27, 5).Value = '  +6.25%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Aptos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(28, 4).Value = '9.77'
$ws.Cells.Item(28, 5).Value = '  +5.94%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +13.10%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '0.997'
$ws.Cells.Item(30, 5).Value = '  -0.40%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.03%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '516.65'
$ws.Cells.Item(32, 5).Value = '  +0.30%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '7.75'
$ws.Cells.Item(33, 5).Value = '  +1.85%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +4.90%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.999'
$ws.Cells.Item(35, 5).Value = '  -0.06%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '164.16'
$ws.Cells.Item(36, 5).Value = '  +0.76%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '19.91'
$ws.Cells.Item(37, 5).Value = '  +4.41%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.93%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '19.38'
$ws.Cells.Item(39, 5).Value = '  +0.30%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '183.84'
$ws.Cells.Item(40, 5).Value = '  +11.67%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.00%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  +5.80%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  +3.96%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '1.69'
$ws.Cells.Item(44, 5).Value = '  +2.49%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +4.09%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '40.01'
$ws.Cells.Item(46, 5).Value = '  +2.47%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '0.0877'
$ws.Cells.Item(47, 5).Value = '  +3.54%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '2.37'
$ws.Cells.Item(48, 5).Value = '  +2.74%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +9.68%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '3.76'
$ws.Cells.Item(50, 5).Value = '  +4.56%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +10.69%  '
